$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a date-looking string into a cell as literal text, without
# letting Excel auto-convert it to a date serial number, and without leaving
# a lingering number-format style on the cell.
function Set-TextValue($cellAddress, $text) {
    $rng = $ws.Range($cellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Update date strings in column A (slash -> dash format) for rows 3-21
Set-TextValue "A3"  "28-07-2022"
Set-TextValue "A4"  "01-08-2022"
Set-TextValue "A5"  "04-08-2022"
Set-TextValue "A6"  "08-08-2022"
Set-TextValue "A7"  "11-08-2022"
Set-TextValue "A8"  "15-08-2022"
Set-TextValue "A9"  "18-08-2022"
Set-TextValue "A10" "22-08-2022"
Set-TextValue "A11" "25-08-2022"
Set-TextValue "A12" "29-08-2022"
Set-TextValue "A13" "01-09-2022"
Set-TextValue "A14" "05-09-2022"
Set-TextValue "A15" "08-09-2022"
Set-TextValue "A16" "12-09-2022"
Set-TextValue "A17" "15-09-2022"
Set-TextValue "A18" "19-09-2022"
Set-TextValue "A19" "22-09-2022"
Set-TextValue "A20" "26-09-2022"
Set-TextValue "A21" "29-09-2022"

# Update attendance-count figures for rows 3-6 (D:Total, E:Real, F:Duplicate, G:Invalid, H:Absent)
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("H6").Value = 0
